# StyleMerger: Use Copy Of baseStyle
#
# Reproduces the issue.4539 edit: turn off gridlines, add three new rows
# (11-13) with a "cellIs equal to" conditional format (highlighting 1 and 2
# with distinct styles), and re-prioritise the existing color-scale rule so
# it evaluates after the new rule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- view: hide gridlines (sheetView/@showGridLines="0") ---------------
$excel.ActiveWindow.DisplayGridlines = $false

# --- new data rows 11:13 -------------------------------------------------
$ws.Range("A11").Value = 1
$ws.Range("A12").Value = 2
$ws.Range("A13").Value = 3

# --- existing colorScale rule on A2:A7: push it behind the new rules ----
$existing = $ws.Range("A2:A7").FormatConditions.Item(1)
$existing.Priority = 3

# --- new conditional formatting on A11:A13 --------------------------------
$fcs = $ws.Range("A11:A13").FormatConditions

# Rule for value = 1 -> thin green (92D050) top border. Ends up dxfId 0.
$ruleOne = $fcs.Add(8, 3, "1")
$ruleOne.Borders.Item(8).LineStyle = 1
$ruleOne.Borders.Item(8).Color = 5296274
$ruleOne.Priority = 1

# Rule for value = 2 -> bold font + thin red (FF0000) top border. Ends up dxfId 1.
$ruleTwo = $fcs.Add(8, 3, "2")
$ruleTwo.Font.Bold = $true
$ruleTwo.Borders.Item(8).LineStyle = 1
$ruleTwo.Borders.Item(8).Color = 255
$ruleTwo.Priority = 2

Write-Output "done"
